$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write cell values in the same order the shared strings were
# originally authored, so the sharedStrings table layout matches.
$ws.Range("A4").Value = "Бахнув Яблок"
$ws.Range("A5").Value = "SCRIPT/G01P06B/c00a0601.ssb"
$ws.Range("A5").Font.Italic = $true
$ws.Range("C4").Value = "Crunch-munch! Chew-chew! Snuffle-slurp!\nCrunch-munch! Chew-chew! Snuffle-slurp!"
$ws.Range("B4").Value = "89 - 104"
$ws.Range("D4").Value = "Хрум-хрум! Ням-ням! Хлюп-хлюп!\nХрум-хрум! Ням-ням! Хлюп-хлюп!"
$ws.Range("E4").Value = "Öñôí-öñôí! Îÿí-îÿí! Öìýð-öìýð!\nÖñôí-öñôí! Îÿí-îÿí! Öìýð-öìýð!"
$ws.Range("C5").Value = " Aaaaah! I\'m done eating! Excellent meal!"
$ws.Range("C6").Value = " I\'m stuffed, and now I\'m getting sleepy…"
$ws.Range("C7").Value = " Yep, I\'m off to bed! Good night, gang!"
$ws.Range("D5").Value = " Ааааах! Я наелся!\nПрекрасный ужин!"
$ws.Range("D6").Value = " Живот набит и теперь я хочу\nспать..."
$ws.Range("D7").Value = " Да, пойду спать!\nДоброй ночи, ребята!"
$ws.Range("E5").Value = " Àààààö! Ÿ îàåìòÿ!\nÐñåëñàòîúê ôçéî!"
$ws.Range("E6").Value = " Çéâïó îàáéó é óåðåñû ÿ öïœô\nòðàóû..."
$ws.Range("E7").Value = " Äà, ðïêäô òðàóû!\nÄïáñïê îïœé, ñåáÿóà!"
$ws.Range("C8").Value = " Good night!"
$ws.Range("D8").Value = " Доброй ночи!"
$ws.Range("E8").Value = " Äïáñïê îïœé!"

# Numeric line-number column
$ws.Range("B5").Value = 115
$ws.Range("B6").Value = 119
$ws.Range("B7").Value = 123
$ws.Range("B8").Value = 127

# Row heights to match source formatting
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(6).RowHeight = 28.8
$ws.Rows.Item(7).RowHeight = 28.8

# Update the view: selection and scroll position per the target sheet view
$ws.Range("D9").Select()
